$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dialog text placeholders from %roomname% to {roomname}
$ws.Range("C2").Value = "I have entered {roomname}"
$ws.Range("C3").Value = "I am in the {roomname}"

# Reflect the final selection left by the editing user
$ws.Range("C3").Select()
